# Apply the "new params to reflect xlsm" update to config.xlsx.
#
# Adds four new config rows (with matching shared-string labels) across the
# General, Demographics and Reporting sheets, tweaks one existing value on
# the Epi sheet, and leaves Demographics as the active/selected sheet -
# mirroring the author's edit.

$wb = $excel.ActiveWorkbook

$general      = $wb.Worksheets.Item("General")
$demographics = $wb.Worksheets.Item("Demographics")
$epi          = $wb.Worksheets.Item("Epi")
$reporting    = $wb.Worksheets.Item("Reporting")

# --- General: two new geography-related parameters ---------------------
$general.Range("A18").Value = "Default_Geography_Initial_Node_Population"
$general.Range("B18").Value = 1000
$general.Range("A19").Value = "Default_Geography_Torus_Size"
$general.Range("B19").Value = 10

# --- Demographics: new minimum adult age parameter ----------------------
$demographics.Range("A19").Value = "Minimum_Adult_Age_Years"
$demographics.Range("B19").Value = 15

# --- Epi: Base_Incubation_Period now enabled (0 -> 1) --------------------
$epi.Range("B4").Value = 1

# --- Reporting: new event-recorder toggle --------------------------------
$reporting.Range("A5").Value = "Report_Event_Recorder"
$reporting.Range("B5").Value = 0

# --- View state: Demographics becomes the active/selected sheet ----------
$general.Range("A20").Select()
$demographics.Activate()
$demographics.Range("A19:B19").Select()
